# Users.xlsx fix: the row for "TestTrainer" (row 11, a leftover test/dummy
# trainer record) was removed. Since the Id column (A) is a simple sequential
# counter independent of the name/contact data, only columns B:F are shifted
# up to close the gap; the Id column itself is left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$srcValues = $ws.Range("B12:F22").Value2
$ws.Range("B11:F21").Value2 = $srcValues
$ws.Rows("22").Delete()
